# ----------------------------------------------------------------------
# digicode.xlsx edit: add 6 new antivirus SKUs (Kaspersky Premium+VPN,
# Norton 360 Deluxe/Premium) as new Tabla1 rows 146-151, tweak a couple
# of existing prices/styles, and move the active selection.
# ----------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Expand Tabla1 (the ListObject) by 6 rows -> A1:K151 ---
$tbl = $ws.ListObjects.Item("Tabla1")
for ($i = 0; $i -lt 6; $i++) {
    [void]$tbl.ListRows.Add()
}

# --- Seed the new rows format from an existing fully-styled data row ---
$ws.Range("A2:K2").Copy()
$ws.Range("A146:K151").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Harmonize the PAY column style for rows 138-145 (drops a stray
#     redundant number-format flag picked up earlier) ---
$ws.Range("F2").Copy()
$ws.Range("F138:F145").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Price corrections on existing McAfee Total Protection rows ---
$ws.Range("H139").Value = 150
$ws.Range("H140").Value = 276

# --- Populate the newly added table rows 146-149 (Kaspersky Premium + VPN) ---
# Row 146
$ws.Range("A146").Value = "DIG10041"
$ws.Range("B146").Value = "Kaspersky Premium + VPN 5PCS"
$ws.Range("C146").Value = 10
$ws.Range("D146").Value = 365
$ws.Range("E146").Value = "Protección PREMIUM"
$ws.Range("F146").Value = "PAY"
$ws.Range("G146").Value = "🤖 Sucripcion x 365 dias."
$ws.Range("H146").Value = 273
$ws.Range("I146").Value = "https://production-tailoy-repo-magento-statics.s3.amazonaws.com/imagenes/872x872/productos/i/a/n/antivirus-kaspersky-premium-10-dispositivos-2-anos-69606-default-1.jpg"

# Row 147
$ws.Range("A147").Value = "DIG10041"
$ws.Range("B147").Value = "Kaspersky Premium + VPN 10PCS"
$ws.Range("C147").Value = 10
$ws.Range("D147").Value = 365
$ws.Range("E147").Value = "Protección PREMIUM"
$ws.Range("F147").Value = "PAY"
$ws.Range("G147").Value = "🤖 Sucripcion x 365 dias."
$ws.Range("H147").Value = 392
$ws.Range("I147").Value = "https://production-tailoy-repo-magento-statics.s3.amazonaws.com/imagenes/872x872/productos/i/a/n/antivirus-kaspersky-premium-10-dispositivos-2-anos-69606-default-1.jpg"

# Row 148
$ws.Range("A148").Value = "DIG10041"
$ws.Range("B148").Value = "Kaspersky Premium + VPN 5PCS"
$ws.Range("C148").Value = 10
$ws.Range("D148").Value = 730
$ws.Range("E148").Value = "Protección PREMIUM"
$ws.Range("F148").Value = "PAY"
$ws.Range("G148").Value = "🤖 Sucripcion x 730 dias."
$ws.Range("H148").Formula = "=+H146*2"
$ws.Range("I148").Value = "https://production-tailoy-repo-magento-statics.s3.amazonaws.com/imagenes/872x872/productos/i/a/n/antivirus-kaspersky-premium-10-dispositivos-2-anos-69606-default-1.jpg"

# Row 149
$ws.Range("A149").Value = "DIG10041"
$ws.Range("B149").Value = "Kaspersky Premium + VPN 10PCS"
$ws.Range("C149").Value = 10
$ws.Range("D149").Value = 730
$ws.Range("E149").Value = "Protección PREMIUM"
$ws.Range("F149").Value = "PAY"
$ws.Range("G149").Value = "🤖 Sucripcion x 730 dias."
$ws.Range("H149").Value = 555
$ws.Range("I149").Value = "https://production-tailoy-repo-magento-statics.s3.amazonaws.com/imagenes/872x872/productos/i/a/n/antivirus-kaspersky-premium-10-dispositivos-2-anos-69606-default-1.jpg"

# --- Populate rows 150-151 (Norton 360) -- names typed first as a pair, ---
# --- then the rest of each row, matching how the author pasted the data ---
$ws.Range("B150").Value = "Norton 360 Deluxe 5PCS"
$ws.Range("B151").Value = "Norton 360 Premium 10PCS"
# Row 150
$ws.Range("A150").Value = "DIG10041"
$ws.Range("C150").Value = 2
$ws.Range("D150").Value = 365
$ws.Range("E150").Value = "ANTIVIRUS"
$ws.Range("F150").Value = "PAY"
$ws.Range("G150").Value = "🤖 Sucripcion x 365 dias."
$ws.Range("H150").Value = 90
$ws.Range("I150").Value = "https://c1.neweggimages.com/productimage/nb1280/32-377-992-12.jpg"

# Row 151
$ws.Range("A151").Value = "DIG10041"
$ws.Range("C151").Value = 2
$ws.Range("D151").Value = 365
$ws.Range("E151").Value = "Protección PREMIUM"
$ws.Range("F151").Value = "PAY"
$ws.Range("G151").Value = "🤖 Sucripcion x 365 dias."
$ws.Range("H151").Value = 120
$ws.Range("I151").Value = "https://c1.neweggimages.com/productimage/nb1280/32-119-004-11.jpg"

# --- Row 151 CODIGO/name cell keeps the bold look used elsewhere in the
#     table (matches the style already on B3) ---
$ws.Range("B3").Copy()
$ws.Range("B151").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B151").Value = "Norton 360 Premium 10PCS"

# --- Hyperlinks for the IMAGEN column on the new rows ---
[void]$ws.Hyperlinks.Add($ws.Range("I146"), "https://production-tailoy-repo-magento-statics.s3.amazonaws.com/imagenes/872x872/productos/i/a/n/antivirus-kaspersky-premium-10-dispositivos-2-anos-69606-default-1.jpg")
[void]$ws.Hyperlinks.Add($ws.Range("I147:I149"), "https://production-tailoy-repo-magento-statics.s3.amazonaws.com/imagenes/872x872/productos/i/a/n/antivirus-kaspersky-premium-10-dispositivos-2-anos-69606-default-1.jpg")
[void]$ws.Hyperlinks.Add($ws.Range("I150"), "https://c1.neweggimages.com/productimage/nb1280/32-377-992-12.jpg")
[void]$ws.Hyperlinks.Add($ws.Range("I151"), "https://c1.neweggimages.com/productimage/nb1280/32-119-004-11.jpg")

# Adding a hyperlink pushes Excel's auto "hyperlink" look onto the cell;
# this sheet's existing IMAGEN hyperlinks keep the plain text style, so
# restore that look on the new cells (relationship itself is untouched).
$ws.Range("I2").Copy()
$ws.Range("I146:I151").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Extend the STOCK "<1" conditional format down through the new rows ---
$ws.Range("C2:C145").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("C2:C151"))

# --- Move the active selection the way the author left the sheet ---
$ws.Range("G152").Select()

Write-Host "Tabla1 range:" $tbl.Range.Address()

